# Generate Report for Handoff
# Adds a new handoff entry (c658e27f-941e-48c1-a98f-0fa0197d0362) as row 3
# on all three worksheets: Overview, zh-cn, de-de.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$commitSha = "204c89b0d2cc013f4b416dd8d009a385cc522dfb"
$baseUrl   = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$commitSha/e2e/"

$fileName = "c658e27f-941e-48c1-a98f-0fa0197d0362.md"
$pathName = "e2e\c658e27f-941e-48c1-a98f-0fa0197d0362.md"

$zhCnXlf = "c658e27f-941e-48c1-a98f-0fa0197d0362.30e7597d9c8b117cc29a1711989144ca3480deed.zh-cn.xlf"
$deDeXlf = "c658e27f-941e-48c1-a98f-0fa0197d0362.30e7597d9c8b117cc29a1711989144ca3480deed.de-de.xlf"

$dateFmt = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------
# Sheet "Overview" (table3.xml) -- expand table A1:G2 -> A1:G3
# ---------------------------------------------------------------------
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A3").Value = $fileName
$wsOverview.Range("B3").Value = $pathName
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), ($baseUrl + $fileName), "", "", $pathName) | Out-Null
$wsOverview.Range("B3").Style = "HyperLink"
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("D3").Value = ""
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-20 06:46:45"
$wsOverview.Range("G3").NumberFormat = $dateFmt

# ---------------------------------------------------------------------
# Sheet "zh-cn" (table1.xml) -- expand table A1:P2 -> A1:P3
# ---------------------------------------------------------------------
$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.ListRows.Add() | Out-Null

$wsZhCn.Range("A3").Value = $fileName
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), ($baseUrl + $fileName), "", "", $fileName) | Out-Null
$wsZhCn.Range("A3").Style = "HyperLink"
$wsZhCn.Range("B3").Value = ".md"
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("D3").Value = "e2e"
$wsZhCn.Range("E3").Value = "ht"
$wsZhCn.Range("F3").Value = "'False"
$wsZhCn.Range("F3").Style = "Normal"
$wsZhCn.Range("G3").Value = $zhCnXlf
$wsZhCn.Range("H3").Value = "2016-08-20 06:46:41"
$wsZhCn.Range("H3").NumberFormat = $dateFmt
$wsZhCn.Range("I3").Value = ""
$wsZhCn.Range("J3").Value = ""
$wsZhCn.Range("K3").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("K3").NumberFormat = $dateFmt
$wsZhCn.Range("L3").Value = ""
$wsZhCn.Range("M3").Value = "'True"
$wsZhCn.Range("M3").Style = "Normal"
$wsZhCn.Range("N3").Value = ""
$wsZhCn.Range("O3").Value = "'False"
$wsZhCn.Range("O3").Style = "Normal"
$wsZhCn.Range("P3").Value = ""

# ---------------------------------------------------------------------
# Sheet "de-de" (table2.xml) -- expand table A1:P2 -> A1:P3
# ---------------------------------------------------------------------
$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.ListRows.Add() | Out-Null

$wsDeDe.Range("A3").Value = $fileName
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), ($baseUrl + $fileName), "", "", $fileName) | Out-Null
$wsDeDe.Range("A3").Style = "HyperLink"
$wsDeDe.Range("B3").Value = ".md"
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("D3").Value = "e2e"
$wsDeDe.Range("E3").Value = "ht"
$wsDeDe.Range("F3").Value = "'False"
$wsDeDe.Range("F3").Style = "Normal"
$wsDeDe.Range("G3").Value = $deDeXlf
$wsDeDe.Range("H3").Value = "2016-08-20 06:46:45"
$wsDeDe.Range("H3").NumberFormat = $dateFmt
$wsDeDe.Range("I3").Value = ""
$wsDeDe.Range("J3").Value = ""
$wsDeDe.Range("K3").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("K3").NumberFormat = $dateFmt
$wsDeDe.Range("L3").Value = ""
$wsDeDe.Range("M3").Value = "'True"
$wsDeDe.Range("M3").Style = "Normal"
$wsDeDe.Range("N3").Value = ""
$wsDeDe.Range("O3").Value = "'False"
$wsDeDe.Range("O3").Style = "Normal"
$wsDeDe.Range("P3").Value = ""

Write-Host "Handoff row added to Overview, zh-cn and de-de sheets."
